$d = $word.ActiveDocument

# Edit 1: "failed to reject the null hypothesis." (split across "fail"/"ed"/" to reject...") -> single run "failed to reject the null hypothesis."
$d.Content.Find.Execute("failed to reject the null hypothesis.", $false, $false, $false, $false, $false, $true, 1, $false, "failed to reject the null hypothesis.", 2)

# Edit 2: "In this random-controlled trial experiment, the consultancy..." -> "Under the randomized controlled trial condition, the consultancy..."
$d.Content.Find.Execute("In this random-controlled trial experiment, the consultancy", $false, $false, $false, $false, $false, $true, 1, $false, "Under the randomized controlled trial condition, the consultancy", 2)

# Edit 3: "In this random-controlled trial experiment, the dosage..." -> "Under the randomized controlled trial condition, the dosage..."
$d.Content.Find.Execute("In this random-controlled trial experiment, the dosage", $false, $false, $false, $false, $false, $true, 1, $false, "Under the randomized controlled trial condition, the dosage", 2)
